$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "80.379.39"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +5.10%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.168.77"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.85%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "208.83"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +4.74%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "624.98"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.38%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.273"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +27.12%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("E9").Value = "  +6.21%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "3.159.26"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.75%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.588"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +26.14%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0000257"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +27.92%  "

$ws.Range("E13").Value = "  +1.48%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "3.744.22"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.27%  "

$ws.Range("E15").Value = "  -0.34%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "31.82"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +7.48%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "80.317.69"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +5.18%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.157.11"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +2.23%  "

$ws.Range("E19").Value = "  +3.52%  "

$ws.Range("E20").Value = "  +9.01%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "9.14"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.48%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "436.21"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +12.66%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.15"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +13.29%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.93"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +7.27%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "3.335.13"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.92%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "75.70"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +3.96%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "4.66"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.80%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "10.83"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +4.25%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("E30").Value = "  +7.71%  "

$ws.Range("E31").Value = "  +0.27%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "8.88"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +5.34%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "558.27"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +9.09%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.45"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.15%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.149"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +12.48%  "

$ws.Range("E36").Value = "  +1.87%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "22.82"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +8.50%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.122"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +18.73%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.27%  "

$ws.Range("E40").Value = "  +5.43%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "20.77"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +3.47%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "162.97"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.47%  "

$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("E44").Value = "  +4.83%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "189.40"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -3.69%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.80"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +6.43%  "

$ws.Range("E47").Value = "  +7.14%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.774"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -3.67%  "

$ws.Range("E49").Value = "  +0.62%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "42.71"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.75%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "4.21"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +6.50%  "

